# Weekly fruit/vegetable price update: add a new week's worth of records
# (4 rows) for Melón / Macroferia Regional de Talca at the top of the
# existing date-ordered block (rows 355-365), pushing the older rows
# down by 4 (they become rows 359-369) without altering their content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows starting at row 355; this shifts the existing
# rows 355:365 down to 359:369 unchanged (values + styles move with them).
$ws.Range("A355:A358").EntireRow.Insert()

# New week's data (Fecha = 44610) for Melón, Calameño/Tuna x Primera/Segunda.
$newRows = @(
  @{ Row=355; H="Calameño"; I="Primera"; J=4000; K=800;  L=800;  M=800;  P=800  },
  @{ Row=356; H="Calameño"; I="Segunda"; J=3000; K=500;  L=500;  M=500;  P=500  },
  @{ Row=357; H="Tuna";     I="Primera"; J=4000; K=800;  L=800;  M=800;  P=800  },
  @{ Row=358; H="Tuna";     I="Segunda"; J=3000; K=500;  L=500;  M=500;  P=500  }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = 5
    $ws.Range("B$row").Value = "Macroferia Regional de Talca"
    $ws.Range("C$row").Value = "Maule"
    $ws.Range("D$row").Value = 44610
    $ws.Range("E$row").Value = 7
    $ws.Range("F$row").Value = 100112027
    $ws.Range("G$row").Value = "Melón"
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $r.I
    $ws.Range("J$row").Value = $r.J
    $ws.Range("K$row").Value = $r.K
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = "$/unidad"
    $ws.Range("O$row").Value = "Región del Maule"
    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = 1
    $ws.Range("R$row").Value = "Hortaliza"
}
